$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.541.50"
$ws.Range("E2").Value = "  +1.53%  "
$ws.Range("D3").Value = "1.878.80"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("E4").Value = "  +0.62%  "
$ws.Range("D5").Value = "'243.50"
$ws.Range("E5").Value = "  +4.94%  "
$ws.Range("D6").Value = "'0.633"
$ws.Range("E6").Value = "  +2.46%  "
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("D8").Value = "'42.84"
$ws.Range("E8").Value = "  +5.19%  "
$ws.Range("D9").Value = "'0.333"
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("D10").Value = "'0.0705"
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("D11").Value = "'0.0994"
$ws.Range("E11").Value = "  +1.54%  "
$ws.Range("D12").Value = "2.147.48"
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'11.85"
$ws.Range("E13").Value = "  +3.35%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.903.60"
$ws.Range("E14").Value = "  +3.16%  "
$ws.Range("D15").Value = "'0.688"
$ws.Range("E15").Value = "  +1.91%  "
$ws.Range("D16").Value = "'4.79"
$ws.Range("E16").Value = "  +2.59%  "
$ws.Range("D17").Value = "35.458.31"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").Value = "'71.22"
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("D19").Value = "0.0₃0807"
$ws.Range("E19").Value = "  +1.95%  "
$ws.Range("D20").Value = "'243.14"
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("D21").Value = "'12.41"
$ws.Range("E21").Value = "  +1.42%  "
$ws.Range("D22").Value = "'4.83"
$ws.Range("E22").Value = "  +3.04%  "
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").Value = "'2.29"
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("D25").Value = "'171.24"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").Value = "'1.99"
$ws.Range("E26").Value = "  +28.52%  "
$ws.Range("D27").Value = "'8.27"
$ws.Range("E27").Value = "  +5.49%  "
$ws.Range("D28").Value = "'17.87"
$ws.Range("E28").Value = "  +1.80%  "
$ws.Range("D29").Value = "'0.125"
$ws.Range("E29").Value = "  +1.28%  "
$ws.Range("D30").Value = "'0.0566"
$ws.Range("E30").Value = "  +1.90%  "
$ws.Range("D31").Value = "'4.08"
$ws.Range("E31").Value = "  +3.48%  "
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").Value = "'4.10"
$ws.Range("E33").Value = "  +3.13%  "
$ws.Range("D34").Value = "'0.897"
$ws.Range("E34").Value = "  +18.22%  "
$ws.Range("D35").Value = "'1.76"
$ws.Range("E35").Value = "  +10.46%  "
$ws.Range("D36").Value = "'2.06"
$ws.Range("E36").Value = "  +4.95%  "
$ws.Range("D37").Value = "'1.36"
$ws.Range("E37").Value = "  +11.41%  "
$ws.Range("D38").Value = "'1.11"
$ws.Range("E38").Value = "  +2.30%  "
$ws.Range("E39").Value = "  +4.15%  "
$ws.Range("D40").Value = "'90.09"
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").Value = "1.360.99"
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("D42").Value = "'15.37"
$ws.Range("E42").Value = "  +5.18%  "
$ws.Range("D43").Value = "'49.30"
$ws.Range("E43").Value = "  +46.00%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'2.37"
$ws.Range("E44").Value = "  +4.76%  "
$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").Value = "'0.0585"
$ws.Range("E45").Value = "  +10.54%  "
$ws.Range("D46").Value = "'12.78"
$ws.Range("E46").Value = "  +47.00%  "
$ws.Range("D47").Value = "'2.43"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("D48").Value = "'6.76"
$ws.Range("E48").Value = "  +7.10%  "
$ws.Range("D49").Value = "'2.73"
$ws.Range("E49").Value = "  -1.44%  "
$ws.Range("D50").Value = "2.063.57"
$ws.Range("E50").Value = "  +1.89%  "
$ws.Range("D51").Value = "'0.0686"
$ws.Range("E51").Value = "  +2.39%  "
